$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, report week dates) ---
$ws.Range("A8").Characters(21, 2).Text = "48"
$ws.Range("C9").Characters(27, 10).Text = "11/24/2025"
$ws.Range("C9").Characters(48, 10).Text = "11/30/2025"

# --- Table data updates (rows 14-33) ---
$ws.Range("L14").Value = -8.333333333333
$ws.Range("N14").Value = -83.582089552238
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "0"
$ws.Range("A15").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("E15").Value = -100
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = -50
$ws.Range("J15").Value = 40
$ws.Range("K15").Value = 30
$ws.Range("L15").Value = 30
$ws.Range("C16").Value = 5
$ws.Range("E16").Value = -44.444444444444
$ws.Range("F16").Value = 35
$ws.Range("G16").Value = 43
$ws.Range("H16").Value = -18.60465116279
$ws.Range("I16").Value = 573
$ws.Range("J16").Value = 586
$ws.Range("K16").Value = -2.218430034129
$ws.Range("L16").Value = -6.065573770491
$ws.Range("M16").Value = 38.072289156626
$ws.Range("N16").Value = -67.387592487194
$ws.Range("C17").Value = 14
$ws.Range("D17").Value = 25
$ws.Range("E17").Value = -44
$ws.Range("F17").Value = 72
$ws.Range("G17").Value = 77
$ws.Range("H17").Value = -6.493506493506
$ws.Range("I17").Value = 957
$ws.Range("J17").Value = 958
$ws.Range("K17").Value = -0.104384133611
$ws.Range("L17").Value = 2.572347266881
$ws.Range("M17").Value = 139.84962406015
$ws.Range("N17").Value = -6.815968841285
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 30
$ws.Range("G18").Value = 35
$ws.Range("H18").Value = -14.285714285714
$ws.Range("I18").Value = 403
$ws.Range("J18").Value = 331
$ws.Range("K18").Value = 21.752265861027
$ws.Range("L18").Value = 37.542662116041
$ws.Range("M18").Value = 116.666666666667
$ws.Range("N18").Value = -67.630522088353
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 14
$ws.Range("E19").Value = -14.285714285714
$ws.Range("F19").Value = 61
$ws.Range("G19").Value = 77
$ws.Range("H19").Value = -20.77922077922
$ws.Range("I19").Value = 893
$ws.Range("J19").Value = 902
$ws.Range("K19").Value = -0.997782705099
$ws.Range("L19").Value = 30.938416422287
$ws.Range("M19").Value = 125.505050505051
$ws.Range("N19").Value = 25.952045133991
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 14
$ws.Range("G20").Value = 19
$ws.Range("H20").Value = -26.315789473684
$ws.Range("I20").Value = 244
$ws.Range("J20").Value = 218
$ws.Range("K20").Value = 11.926605504587
$ws.Range("L20").Value = -20.521172638436
$ws.Range("M20").Value = 117.857142857143
$ws.Range("N20").Value = -60.325203252032
$ws.Range("C21").Value = 37
$ws.Range("D21").Value = 55
$ws.Range("E21").Value = -32.727272727272
$ws.Range("F21").Value = 213
$ws.Range("G21").Value = 255
$ws.Range("H21").Value = -16.470588235294
$ws.Range("I21").Value = 3133
$ws.Range("J21").Value = 3049
$ws.Range("K21").Value = 2.755001639881
$ws.Range("L21").Value = 8.898157803267
$ws.Range("M21").Value = 102.914507772021
$ws.Range("N21").Value = -42.974153622133
$ws.Range("C22").Value = 2
$ws.Range("D22").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 4
$ws.Range("G22").Value = 10
$ws.Range("H22").Value = -60
$ws.Range("I22").Value = 64
$ws.Range("J22").Value = 84
$ws.Range("K22").Value = -23.809523809523
$ws.Range("L22").Value = 3.225806451612
$ws.Range("M22").Value = 12.280701754386
$ws.Range("C23").Value = 7
$ws.Range("D23").Value = 5
$ws.Range("E23").Value = 40
$ws.Range("F23").Value = 30
$ws.Range("G23").Value = 33
$ws.Range("H23").Value = -9.090909090909
$ws.Range("I23").Value = 463
$ws.Range("J23").Value = 450
$ws.Range("K23").Value = 2.888888888888
$ws.Range("L23").Value = 3.811659192825
$ws.Range("M23").Value = 67.148014440433
$ws.Range("C24").Value = 29
$ws.Range("D24").Value = 32
$ws.Range("E24").Value = -9.375
$ws.Range("F24").Value = 139
$ws.Range("G24").Value = 131
$ws.Range("H24").Value = 6.106870229007
$ws.Range("I24").Value = 1837
$ws.Range("J24").Value = 1574
$ws.Range("K24").Value = 16.709021601016
$ws.Range("L24").Value = 23.703703703703
$ws.Range("M24").Value = 41.853281853281
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 16
$ws.Range("E25").Value = -62.5
$ws.Range("F25").Value = 40
$ws.Range("G25").Value = 62
$ws.Range("H25").Value = -35.483870967741
$ws.Range("I25").Value = 747
$ws.Range("J25").Value = 694
$ws.Range("K25").Value = 7.636887608069
$ws.Range("L25").Value = 30.823117338003
$ws.Range("C26").Value = 30
$ws.Range("D26").Value = 23
$ws.Range("E26").Value = 30.434782608695
$ws.Range("F26").Value = 100
$ws.Range("G26").Value = 92
$ws.Range("H26").Value = 8.695652173913
$ws.Range("I26").Value = 1190
$ws.Range("J26").Value = 1159
$ws.Range("K26").Value = 2.674719585849
$ws.Range("L26").Value = 13.875598086124
$ws.Range("M26").Value = 11.214953271028
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "0"
$ws.Range("A27").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("E27").Value = -100
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = -50
$ws.Range("J27").Value = 52
$ws.Range("K27").Value = 23.076923076923
$ws.Range("L27").Value = 10.344827586206
$ws.Range("F28").Value = 14
$ws.Range("G28").Value = 8
$ws.Range("H28").Value = 75
$ws.Range("I28").Value = 113
$ws.Range("J28").Value = 150
$ws.Range("K28").Value = -24.666666666666
$ws.Range("L28").Value = 10.78431372549
$ws.Range("D29").Value = 1
$ws.Range("J29").Value = 51
$ws.Range("K29").Value = -35.294117647058
$ws.Range("L29").Value = -2.941176470588
$ws.Range("N29").Value = -83.823529411764
$ws.Range("D30").Value = 1
$ws.Range("J30").Value = 46
$ws.Range("K30").Value = -41.304347826087
$ws.Range("L30").Value = -15.625
$ws.Range("N30").Value = -85.326086956521
$ws.Range("F33").NumberFormat = "@"
$ws.Range("F33").Value = "0"
$ws.Range("A33").Copy()
$ws.Range("F33").PasteSpecial(-4122)

$excel.CutCopyMode = $false
